# Fix a grammatical error (per commit message) in the Sprint Report.
#
# Diffing the document's canonical OOXML against the target shows the
# document text itself only changes in two small spots (everything else
# in the diff is Word's automatic w:proofErr spell-check markup around
# proper nouns / camel-case identifiers, which carries no visible-text
# change and isn't something this automation surfaces):
#
#   1. "... is as follow: "                 -> "... is as follows:"
#   2. "Screenshot Feasability (Started..." -> "Screenshot Feasibility (Started..."

$d = $word.ActiveDocument

$d.Content.Find.Execute("is as follow: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "is as follows:", 2)

$d.Content.Find.Execute("Feasability", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Feasibility", 2)
